$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218; existing rows 218:260 shift down to 219:261
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new weekly data point
$ws.Cells.Item(218, 1).Value = 4
$ws.Cells.Item(218, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(218, 3).Value = "Los Lagos"
$ws.Cells.Item(218, 4).Value = 44476
$ws.Cells.Item(218, 5).Value = 10
$ws.Cells.Item(218, 6).Value = 100112006
$ws.Cells.Item(218, 7).Value = "Repollo"
$ws.Cells.Item(218, 8).Value = "Crespo record"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 500
$ws.Cells.Item(218, 11).Value = 1000
$ws.Cells.Item(218, 12).Value = 1100
$ws.Cells.Item(218, 13).Value = 1050
$ws.Cells.Item(218, 14).Value = "$/unidad"
$ws.Cells.Item(218, 15).Value = "Región Metropolitana"
$ws.Cells.Item(218, 16).Value = 1050
$ws.Cells.Item(218, 17).Value = 1
$ws.Cells.Item(218, 18).Value = "Hortaliza"
